# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the per-item profit tables on each sheet,
# matching the scheduled-runner refresh of computed market-board columns
# (currentAveragePrice*, LevePrice*, LeveProfit*).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 35000
$ws.Range("I47").Value = 35000
$ws.Range("K47").Value = 35000
$ws.Range("M47").Value = -34028
# Row 137
$ws.Range("H137").Value = 5665.4
$ws.Range("I137").Value = 4442.3335
$ws.Range("K137").Value = 13327.0005
$ws.Range("M137").Value = -10777.0005
# Row 138
$ws.Range("H138").Value = 6451.037
$ws.Range("I138").Value = 3916.6667
$ws.Range("J138").Value = 6767.8335
$ws.Range("K138").Value = 11750.0001
$ws.Range("L138").Value = 20303.5005
$ws.Range("M138").Value = -6610.000100000001
$ws.Range("N138").Value = -30583.5005

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23730.55
$ws.Range("I32").Value = 2488.5
$ws.Range("K32").Value = 2488.5
$ws.Range("M32").Value = -2201.5
# Row 38
$ws.Range("H38").Value = 8518
$ws.Range("I38").Value = 8518
$ws.Range("K38").Value = 8518
$ws.Range("M38").Value = -8051
# Row 74
$ws.Range("H74").Value = 2283.3489
$ws.Range("I74").Value = 2282.2307
$ws.Range("K74").Value = 2282.2307
$ws.Range("M74").Value = -1408.2307
# Row 77
$ws.Range("H77").Value = 2283.3489
$ws.Range("I77").Value = 2282.2307
$ws.Range("K77").Value = 11411.1535
$ws.Range("M77").Value = -7043.1535

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4373.55
$ws.Range("I86").Value = 2989.7693
$ws.Range("K86").Value = 2989.7693
$ws.Range("M86").Value = -1866.7693
# Row 89
$ws.Range("H89").Value = 4373.55
$ws.Range("I89").Value = 2989.7693
$ws.Range("K89").Value = 14948.8465
$ws.Range("M89").Value = -9332.8465
# Row 94
$ws.Range("H94").Value = 1968.0294
$ws.Range("I94").Value = 772.7143
$ws.Range("J94").Value = 3898.923
$ws.Range("K94").Value = 772.7143
$ws.Range("L94").Value = 3898.923
$ws.Range("M94").Value = -321.7143
$ws.Range("N94").Value = -4800.923
# Row 105
$ws.Range("H105").Value = 7786.6665
$ws.Range("I105").Value = 7869.7
$ws.Range("J105").Value = 7371.5
$ws.Range("K105").Value = 7869.7
$ws.Range("L105").Value = 7371.5
$ws.Range("M105").Value = -6122.7
$ws.Range("N105").Value = -10865.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 82.833336
$ws.Range("I7").Value = 75
$ws.Range("K7").Value = 75
$ws.Range("M7").Value = 38
# Row 31
$ws.Range("H31").Value = 4821.738
$ws.Range("I31").Value = 4900.44
$ws.Range("J31").Value = 4464
$ws.Range("K31").Value = 4900.44
$ws.Range("L31").Value = 4464
$ws.Range("M31").Value = -4605.44
$ws.Range("N31").Value = -5054
# Row 34
$ws.Range("H34").Value = 4821.738
$ws.Range("I34").Value = 4900.44
$ws.Range("J34").Value = 4464
$ws.Range("K34").Value = 4900.44
$ws.Range("L34").Value = 4464
$ws.Range("M34").Value = -4698.44
$ws.Range("N34").Value = -4868
# Row 47
$ws.Range("H47").Value = 34999.5
$ws.Range("I47").Value = 30000
$ws.Range("J47").Value = 39999
$ws.Range("K47").Value = 30000
$ws.Range("L47").Value = 39999
$ws.Range("M47").Value = -29434
$ws.Range("N47").Value = -41131
# Row 58
$ws.Range("H58").Value = 5229.1904
$ws.Range("I58").Value = 5636.0586
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 5636.0586
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -5433.0586
$ws.Range("N58").Value = -3906
# Row 134
$ws.Range("H134").Value = 2818.5518
$ws.Range("I134").Value = 1057.4762
$ws.Range("J134").Value = 7441.375
$ws.Range("K134").Value = 3172.4286
$ws.Range("L134").Value = 22324.125
$ws.Range("M134").Value = -637.4286000000002
$ws.Range("N134").Value = -27394.125
# Row 136
$ws.Range("H136").Value = 5229.1904
$ws.Range("I136").Value = 5636.0586
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 16908.1758
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -14358.1758
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 803.55554
$ws.Range("I5").Value = 664.1667
$ws.Range("J5").Value = 1082.3334
$ws.Range("K5").Value = 1992.5001
$ws.Range("L5").Value = 3247.0002
$ws.Range("M5").Value = -1880.5001
$ws.Range("N5").Value = -3471.0002
# Row 135
$ws.Range("H135").Value = 803.55554
$ws.Range("I135").Value = 664.1667
$ws.Range("J135").Value = 1082.3334
$ws.Range("K135").Value = 5977.5003
$ws.Range("L135").Value = 9741.000599999999
$ws.Range("M135").Value = -3442.5003
$ws.Range("N135").Value = -14811.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3405.9333
$ws.Range("I80").Value = 2703.9
$ws.Range("J80").Value = 4810
$ws.Range("K80").Value = 2703.9
$ws.Range("L80").Value = 4810
$ws.Range("M80").Value = -1705.9
$ws.Range("N80").Value = -6806
# Row 83
$ws.Range("H83").Value = 3405.9333
$ws.Range("I83").Value = 2703.9
$ws.Range("J83").Value = 4810
$ws.Range("K83").Value = 13519.5
$ws.Range("L83").Value = 24050
$ws.Range("M83").Value = -8527.5
$ws.Range("N83").Value = -34034
# Row 97
$ws.Range("H97").Value = 216.5
$ws.Range("I97").Value = 216.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 216.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 279.5
$ws.Range("N97").ClearContents()
# Row 126
$ws.Range("H126").Value = 7447.95
$ws.Range("I126").Value = 7153.75
$ws.Range("K126").Value = 21461.25
$ws.Range("M126").Value = -18991.25
# Row 132
$ws.Range("H132").Value = 10957.667
$ws.Range("I132").Value = 10957.667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 32873.001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -30343.001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1872.7778
$ws.Range("I16").Value = 1981.625
$ws.Range("J16").Value = 1002
$ws.Range("K16").Value = 1981.625
$ws.Range("L16").Value = 1002
$ws.Range("M16").Value = -1811.625
$ws.Range("N16").Value = -1342
# Row 22
$ws.Range("H22").Value = 2211
$ws.Range("I22").Value = 1057.8
$ws.Range("J22").Value = 3034.7144
$ws.Range("K22").Value = 1057.8
$ws.Range("L22").Value = 3034.7144
$ws.Range("M22").Value = -762.8
$ws.Range("N22").Value = -3624.7144
# Row 27
$ws.Range("H27").Value = 2211
$ws.Range("I27").Value = 1057.8
$ws.Range("J27").Value = 3034.7144
$ws.Range("K27").Value = 1057.8
$ws.Range("L27").Value = 3034.7144
$ws.Range("M27").Value = -950.8
$ws.Range("N27").Value = -3248.7144
# Row 100
$ws.Range("H100").Value = 4877.625
$ws.Range("I100").Value = 3049.2856
$ws.Range("J100").Value = 6299.6665
$ws.Range("K100").Value = 3049.2856
$ws.Range("L100").Value = 6299.6665
$ws.Range("M100").Value = -2508.2856
$ws.Range("N100").Value = -7381.6665

Write-Output "Updated market-board columns across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets."